# edit.ps1 - PowerPoint COM-interop script
#
# Applies the changes described by the target diff:
#   1. Updates the cached "datetimeFigureOut" date field text from
#      6/11/2024 -> 6/13/2024 on the slide master and on every slide
#      layout that carries that placeholder.
#   2. Appends a new slide 11 (using the "Title Slide" layout) with a
#      title of "Project Link :" and a subtitle containing a hyperlinked
#      URL to the project's GitHub repository.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh the cached date placeholder text everywhere it appears.
# ---------------------------------------------------------------------
$oldDate = "6/11/2024"
$newDate = "6/13/2024"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Add the new "Project Link" slide at the end of the deck.
# ---------------------------------------------------------------------
$titleSlideLayout = $layouts.Item(1)   # "Title Slide" custom layout
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.AddSlide($newIndex, $titleSlideLayout)

$titleShape = $s.Shapes.Item(1)
$titleShape.Name = "Title 1"
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Text = "Project Link :"
$titleTr.LanguageID = "en-IN"

$subtitleShape = $s.Shapes.Item(2)
$subtitleShape.Name = "Subtitle 2"
$subTr = $subtitleShape.TextFrame.TextRange
$subTr.Text = "https://github.com/LokeshSai29/key_logger_project_APSSDC-CS.git"
$subTr.LanguageID = "en-IN"

$action = $subTr.ActionSettings.Item(1)
$action.Hyperlink.Address = "https://github.com/LokeshSai29/key_logger_project_APSSDC-CS.git"

Write-Output "edit.ps1 completed"
